$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A4:M7").ClearContents()
$ws.Range("J11").Select()
